# Insert a new data row at row 498 (new price observation for Ajo Chino,
# dated 2023-10-13 / serial 45212), shifting the existing rows 498-554
# down to 499-555, exactly as in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("498").Insert()

$ws.Range("A498").Value = 8
$ws.Range("B498").Value = "Terminal La Palmera de La Serena"
$ws.Range("C498").Value = "Coquimbo"
$ws.Range("D498").Value = 45212
$ws.Range("E498").Value = 4
$ws.Range("F498").Value = 100112003
$ws.Range("G498").Value = "Ajo"
$ws.Range("H498").Value = "Chino"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 450
$ws.Range("K498").Value = 24000
$ws.Range("L498").Value = 25000
$ws.Range("M498").Value = 24500
$ws.Range("N498").Value = "$/malla 10 kilos"
$ws.Range("O498").Value = "China"
$ws.Range("P498").Value = 2450
$ws.Range("Q498").Value = 10
$ws.Range("R498").Value = "Hortaliza"
